$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 618.0645
$ws.Range("I28").Value = 249.61111
$ws.Range("J28").Value = 1128.2307
$ws.Range("K28").Value = 249.61111
$ws.Range("L28").Value = 1128.2307
$ws.Range("M28").Value = 235.38889
$ws.Range("N28").Value = -2098.2307

$ws.Range("H92").Value = 1437.375
$ws.Range("I92").Value = 1600
$ws.Range("J92").Value = 949.5
$ws.Range("K92").Value = 1600
$ws.Range("L92").Value = 949.5
$ws.Range("M92").Value = -352
$ws.Range("N92").Value = -3445.5

$ws.Range("H96").Value = 791.5
$ws.Range("I96").Value = 501.5
$ws.Range("J96").Value = 984.8333
$ws.Range("K96").Value = 1504.5
$ws.Range("L96").Value = 2954.4999
$ws.Range("M96").Value = -131.5
$ws.Range("N96").Value = -5700.4999

$ws.Range("H98").Value = 1696.4615
$ws.Range("I98").Value = 1680
$ws.Range("J98").Value = 1894
$ws.Range("K98").Value = 1680
$ws.Range("L98").Value = 1894
$ws.Range("M98").Value = -182
$ws.Range("N98").Value = -4890

$ws.Range("H113").Value = 102661
$ws.Range("I113").Value = 145801.42
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 145801.42
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -142547.42
$ws.Range("N113").Value = -8508

$ws.Range("H116").Value = 1884.8889
$ws.Range("I116").Value = 1851.4286
$ws.Range("J116").Value = 2002
$ws.Range("K116").Value = 1851.4286
$ws.Range("L116").Value = 2002
$ws.Range("M116").Value = 1590.5714
$ws.Range("N116").Value = -8886

$ws.Range("H122").Value = 1696.4615
$ws.Range("I122").Value = 1680
$ws.Range("J122").Value = 1894
$ws.Range("K122").Value = 5040
$ws.Range("L122").Value = 5682
$ws.Range("M122").Value = -2590
$ws.Range("N122").Value = -10582

$ws.Range("H132").Value = 5004541
$ws.Range("I132").Value = 5560534
$ws.Range("J132").Value = 603.6
$ws.Range("K132").Value = 16681602
$ws.Range("L132").Value = 1810.8
$ws.Range("M132").Value = -16679072
$ws.Range("N132").Value = -6870.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 336493.34
$ws.Range("I2").Value = 5980
$ws.Range("J2").Value = 501750
$ws.Range("K2").Value = 5980
$ws.Range("L2").Value = 501750
$ws.Range("M2").Value = -5867
$ws.Range("N2").Value = -501976

$ws.Range("H24").Value = 25177.5
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 25177.5
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 25177.5
$ws.Range("N24").Value = -25925.5

$ws.Range("H32").Value = 22138.11
$ws.Range("I32").Value = 3822.125
$ws.Range("J32").Value = 142024.55
$ws.Range("K32").Value = 3822.125
$ws.Range("L32").Value = 142024.55
$ws.Range("M32").Value = -3535.125

$ws.Range("H74").Value = 3227315.5
$ws.Range("I74").Value = 1066.8422
$ws.Range("J74").Value = 8335543
$ws.Range("K74").Value = 1066.8422
$ws.Range("L74").Value = 8335543
$ws.Range("M74").Value = -192.8422
$ws.Range("N74").Value = -8337291

$ws.Range("H77").Value = 3227315.5
$ws.Range("I77").Value = 1066.8422
$ws.Range("J77").Value = 8335543
$ws.Range("K77").Value = 5334.211
$ws.Range("L77").Value = 41677715
$ws.Range("M77").Value = -966.2110000000002
$ws.Range("N77").Value = -41686451

$ws.Range("H100").Value = 25177.5
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 25177.5
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 25177.5
$ws.Range("N100").Value = -27341.5

$ws.Range("H116").Value = 336493.34
$ws.Range("I116").Value = 5980
$ws.Range("J116").Value = 501750
$ws.Range("K116").Value = 5980
$ws.Range("L116").Value = 501750
$ws.Range("M116").Value = -3686
$ws.Range("N116").Value = -506338

$ws.Range("H122").Value = 1589.9286
$ws.Range("I122").Value = 1531.5834
$ws.Range("J122").Value = 1940
$ws.Range("K122").Value = 4594.7502
$ws.Range("L122").Value = 5820
$ws.Range("M122").Value = -2144.7502
$ws.Range("N122").Value = -10720

$ws.Range("H133").Value = 60000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 60000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 60000
$ws.Range("N133").Value = -65060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 336493.34
$ws.Range("I3").Value = 5980
$ws.Range("J3").Value = 501750
$ws.Range("K3").Value = 5980
$ws.Range("L3").Value = 501750
$ws.Range("M3").Value = -5866
$ws.Range("N3").Value = -501978

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H118").Value = 46303
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 46303
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 46303
$ws.Range("N118").Value = -49617

$ws.Range("H122").Value = 537.3333
$ws.Range("I122").Value = 537.3333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1611.9999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 838.0001
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 3593.0244
$ws.Range("I132").Value = 3487.5625
$ws.Range("J132").Value = 3968
$ws.Range("K132").Value = 10462.6875
$ws.Range("L132").Value = 11904
$ws.Range("M132").Value = -7932.6875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 42.5
$ws.Range("I7").Value = 45
$ws.Range("J7").Value = 40
$ws.Range("K7").Value = 135
$ws.Range("L7").Value = 120
$ws.Range("M7").Value = -23
$ws.Range("N7").Value = -344

$ws.Range("H114").Value = 494.91666
$ws.Range("I114").Value = 264.36365
$ws.Range("J114").Value = 3031
$ws.Range("K114").Value = 793.09095
$ws.Range("L114").Value = 9093
$ws.Range("M114").Value = 2460.90905

$ws.Range("H117").Value = 11115.909
$ws.Range("I117").Value = 362.2
$ws.Range("J117").Value = 20077.334
$ws.Range("K117").Value = 1086.6
$ws.Range("L117").Value = 60232.00199999999
$ws.Range("M117").Value = 2355.4
$ws.Range("N117").Value = -67116.00199999999

$ws.Range("H121").Value = 2996.818
$ws.Range("I121").Value = 1590
$ws.Range("J121").Value = 3970.7693
$ws.Range("K121").Value = 4770
$ws.Range("L121").Value = 11912.3079
$ws.Range("M121").Value = -3460

$ws.Range("H122").Value = 606.2
$ws.Range("I122").Value = 482.1111
$ws.Range("J122").Value = 792.3333
$ws.Range("K122").Value = 4338.9999
$ws.Range("L122").Value = 7130.9997
$ws.Range("M122").Value = -1888.9999
$ws.Range("N122").Value = -12030.9997

$ws.Range("H131").Value = 807.9798
$ws.Range("I131").Value = 464.63635
$ws.Range("J131").Value = 850.8977
$ws.Range("K131").Value = 1393.90905
$ws.Range("L131").Value = 2552.6931
$ws.Range("M131").Value = 3646.09095
$ws.Range("N131").Value = -12632.6931

$ws.Range("H132").Value = 2045
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 2343.5715
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 21092.1435
$ws.Range("M132").Value = -6470
$ws.Range("N132").Value = -26152.1435

$ws.Range("H134").Value = 2683.3333
$ws.Range("I134").Value = 2281.5386
$ws.Range("J134").Value = 2944.5
$ws.Range("K134").Value = 6844.6158
$ws.Range("L134").Value = 8833.5
$ws.Range("M134").Value = -1774.6158
$ws.Range("N134").Value = -18973.5

$ws.Range("H140").Value = 5850.391
$ws.Range("I140").Value = 7777.2666
$ws.Range("J140").Value = 2237.5
$ws.Range("K140").Value = 23331.7998
$ws.Range("L140").Value = 6712.5
$ws.Range("M140").Value = -18151.7998
$ws.Range("N140").Value = -17072.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws.Range("H135").Value = 28586
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 28586
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 28586
$ws.Range("N135").Value = -38726

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5164
$ws.Range("I61").Value = 4334.6665
$ws.Range("J61").Value = 5993.3335
$ws.Range("K61").Value = 4334.6665
$ws.Range("L61").Value = 5993.3335
$ws.Range("M61").Value = -4132.6665

$ws.Range("H98").Value = 26000
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 26000
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 26000
$ws.Range("N98").Value = -31990

$ws.Range("H113").Value = 5164
$ws.Range("I113").Value = 4334.6665
$ws.Range("J113").Value = 5993.3335
$ws.Range("K113").Value = 4334.6665
$ws.Range("L113").Value = 5993.3335
$ws.Range("M113").Value = -2164.6665

$ws.Range("H122").Value = 1600.5714

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 1000000000
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 1000000000
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 1000000000
$ws.Range("N95").Value = -1000005492

$ws.Range("H122").Value = 1500
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -10900

$ws.Range("H136").Value = 771.381
$ws.Range("I136").Value = 636
$ws.Range("J136").Value = 1042.1428
$ws.Range("K136").Value = 1908
$ws.Range("L136").Value = 3126.4284
$ws.Range("M136").Value = 642
$ws.Range("N136").Value = -8226.428400000001
